$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-11-08, serial 44508) is inserted as row 25,
# pushing the existing rows 25-32 down to rows 26-33.
$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25, 1).Value = 3
$ws.Cells.Item(25, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 44508
$ws.Cells.Item(25, 5).Value = 5
$ws.Cells.Item(25, 6).Value = 100112022
$ws.Cells.Item(25, 7).Value = "Arveja Verde"
$ws.Cells.Item(25, 8).Value = "Perfection"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 68
$ws.Cells.Item(25, 11).Value = 16000
$ws.Cells.Item(25, 12).Value = 17000
$ws.Cells.Item(25, 13).Value = 16515
$ws.Cells.Item(25, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 661
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
